$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 22.614608
$ws.Cells.Item(2, 8).Value = 67.843824
$ws.Cells.Item(2, 9).Value = 0.3650188533124966
$ws.Cells.Item(2, 10).Value = 0.3650188533124966
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 10.39091433333333
$ws.Cells.Item(2, 14).Value = 31.172743
$ws.Cells.Item(2, 15).Value = 0.0835098648954196
$ws.Cells.Item(2, 16).Value = 0.0835098648954196
$ws.Cells.Item(2, 17).Value = 234.9864544099147
$ws.Cells.Item(2, 18).Value = 2114.878089689232
$ws.Cells.Item(2, 19).Value = 0.03048267512440758
$ws.Cells.Item(2, 20).Value = 0.03048267512440758

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 22.614608
$ws.Cells.Item(3, 8).Value = 67.843824
$ws.Cells.Item(3, 9).Value = 0.3650188533124966
$ws.Cells.Item(3, 10).Value = 0.3650188533124966
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 35.71561933333334
$ws.Cells.Item(3, 14).Value = 107.146858
$ws.Cells.Item(3, 15).Value = 0.287039855156433
$ws.Cells.Item(3, 16).Value = 0.287039855156433
$ws.Cells.Item(3, 17).Value = 807.6947307005547
$ws.Cells.Item(3, 18).Value = 7269.252576304993
$ws.Cells.Item(3, 19).Value = 0.1047749587841863
$ws.Cells.Item(3, 20).Value = 0.1047749587841863

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 22.614608
$ws.Cells.Item(4, 8).Value = 67.843824
$ws.Cells.Item(4, 9).Value = 0.3650188533124966
$ws.Cells.Item(4, 10).Value = 0.3650188533124966
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 51.87044599999999
$ws.Cells.Item(4, 14).Value = 155.611338
$ws.Cells.Item(4, 15).Value = 0.4168732219867682
$ws.Cells.Item(4, 16).Value = 0.4168732219867682
$ws.Cells.Item(4, 17).Value = 1173.029803075168
$ws.Cells.Item(4, 18).Value = 10557.26822767651
$ws.Cells.Item(4, 19).Value = 0.152166585466296
$ws.Cells.Item(4, 20).Value = 0.152166585466296

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 22.614608
$ws.Cells.Item(5, 8).Value = 67.843824
$ws.Cells.Item(5, 9).Value = 0.3650188533124966
$ws.Cells.Item(5, 10).Value = 0.3650188533124966
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 26.450408
$ws.Cells.Item(5, 14).Value = 79.351224
$ws.Cells.Item(5, 15).Value = 0.2125770579613792
$ws.Cells.Item(5, 16).Value = 0.2125770579613792
$ws.Cells.Item(5, 17).Value = 598.165608360064
$ws.Cells.Item(5, 18).Value = 5383.490475240576
$ws.Cells.Item(5, 19).Value = 0.07759463393760677
$ws.Cells.Item(5, 20).Value = 0.07759463393760677

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 12.38193366666667
$ws.Cells.Item(6, 8).Value = 37.145801
$ws.Cells.Item(6, 9).Value = 0.1998548561530699
$ws.Cells.Item(6, 10).Value = 0.1998548561530699
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 10.39091433333333
$ws.Cells.Item(6, 14).Value = 31.172743
$ws.Cells.Item(6, 15).Value = 0.0835098648954196
$ws.Cells.Item(6, 16).Value = 0.0835098648954196
$ws.Cells.Item(6, 17).Value = 128.6596120113492
$ws.Cells.Item(6, 18).Value = 1157.936508102143
$ws.Cells.Item(6, 19).Value = 0.01668985203603638
$ws.Cells.Item(6, 20).Value = 0.01668985203603638

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 12.38193366666667
$ws.Cells.Item(7, 8).Value = 37.145801
$ws.Cells.Item(7, 9).Value = 0.1998548561530699
$ws.Cells.Item(7, 10).Value = 0.1998548561530699
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 35.71561933333334
$ws.Cells.Item(7, 14).Value = 107.146858
$ws.Cells.Item(7, 15).Value = 0.287039855156433
$ws.Cells.Item(7, 16).Value = 0.287039855156433
$ws.Cells.Item(7, 17).Value = 442.2284294492509
$ws.Cells.Item(7, 18).Value = 3980.055865043258
$ws.Cells.Item(7, 19).Value = 0.05736630896248693
$ws.Cells.Item(7, 20).Value = 0.05736630896248692

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 12.38193366666667
$ws.Cells.Item(8, 8).Value = 37.145801
$ws.Cells.Item(8, 9).Value = 0.1998548561530699
$ws.Cells.Item(8, 10).Value = 0.1998548561530699
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 51.87044599999999
$ws.Cells.Item(8, 14).Value = 155.611338
$ws.Cells.Item(8, 15).Value = 0.4168732219867682
$ws.Cells.Item(8, 16).Value = 0.4168732219867682
$ws.Cells.Item(8, 17).Value = 642.2564216324153
$ws.Cells.Item(8, 18).Value = 5780.307794691737
$ws.Cells.Item(8, 19).Value = 0.08331413781423233
$ws.Cells.Item(8, 20).Value = 0.08331413781423233

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 12.38193366666667
$ws.Cells.Item(9, 8).Value = 37.145801
$ws.Cells.Item(9, 9).Value = 0.1998548561530699
$ws.Cells.Item(9, 10).Value = 0.1998548561530699
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 26.450408
$ws.Cells.Item(9, 14).Value = 79.351224
$ws.Cells.Item(9, 15).Value = 0.2125770579613792
$ws.Cells.Item(9, 16).Value = 0.2125770579613792
$ws.Cells.Item(9, 17).Value = 327.5071973122693
$ws.Cells.Item(9, 18).Value = 2947.564775810424
$ws.Cells.Item(9, 19).Value = 0.04248455734031424
$ws.Cells.Item(9, 20).Value = 0.04248455734031424

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 19.82277733333333
$ws.Cells.Item(10, 8).Value = 59.468332
$ws.Cells.Item(10, 9).Value = 0.3199563508543806
$ws.Cells.Item(10, 10).Value = 0.3199563508543806
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 10.39091433333333
$ws.Cells.Item(10, 14).Value = 31.172743
$ws.Cells.Item(10, 15).Value = 0.0835098648954196
$ws.Cells.Item(10, 16).Value = 0.0835098648954196
$ws.Cells.Item(10, 17).Value = 205.9767811194084
$ws.Cells.Item(10, 18).Value = 1853.791030074676
$ws.Cells.Item(10, 19).Value = 0.0267195116322808
$ws.Cells.Item(10, 20).Value = 0.0267195116322808

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 19.82277733333333
$ws.Cells.Item(11, 8).Value = 59.468332
$ws.Cells.Item(11, 9).Value = 0.3199563508543806
$ws.Cells.Item(11, 10).Value = 0.3199563508543806
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 35.71561933333334
$ws.Cells.Item(11, 14).Value = 107.146858
$ws.Cells.Item(11, 15).Value = 0.287039855156433
$ws.Cells.Item(11, 16).Value = 0.287039855156433
$ws.Cells.Item(11, 17).Value = 707.9827693667618
$ws.Cells.Item(11, 18).Value = 6371.844924300856
$ws.Cells.Item(11, 19).Value = 0.09184022460562226
$ws.Cells.Item(11, 20).Value = 0.09184022460562226

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 19.82277733333333
$ws.Cells.Item(12, 8).Value = 59.468332
$ws.Cells.Item(12, 9).Value = 0.3199563508543806
$ws.Cells.Item(12, 10).Value = 0.3199563508543806
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 51.87044599999999
$ws.Cells.Item(12, 14).Value = 155.611338
$ws.Cells.Item(12, 15).Value = 0.4168732219867682
$ws.Cells.Item(12, 16).Value = 0.4168732219867682
$ws.Cells.Item(12, 17).Value = 1028.21630123869
$ws.Cells.Item(12, 18).Value = 9253.946711148215
$ws.Cells.Item(12, 19).Value = 0.1333812348757945
$ws.Cells.Item(12, 20).Value = 0.1333812348757945

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 19.82277733333333
$ws.Cells.Item(13, 8).Value = 59.468332
$ws.Cells.Item(13, 9).Value = 0.3199563508543806
$ws.Cells.Item(13, 10).Value = 0.3199563508543806
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 26.450408
$ws.Cells.Item(13, 14).Value = 79.351224
$ws.Cells.Item(13, 15).Value = 0.2125770579613792
$ws.Cells.Item(13, 16).Value = 0.2125770579613792
$ws.Cells.Item(13, 17).Value = 524.3205481598186
$ws.Cells.Item(13, 18).Value = 4718.884933438368
$ws.Cells.Item(13, 19).Value = 0.06801537974068304
$ws.Cells.Item(13, 20).Value = 0.06801537974068304

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 7.135311000000001
$ws.Cells.Item(14, 8).Value = 21.405933
$ws.Cells.Item(14, 9).Value = 0.115169939680053
$ws.Cells.Item(14, 10).Value = 0.115169939680053
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 10.39091433333333
$ws.Cells.Item(14, 14).Value = 31.172743
$ws.Cells.Item(14, 15).Value = 0.0835098648954196
$ws.Cells.Item(14, 16).Value = 0.0835098648954196
$ws.Cells.Item(14, 17).Value = 74.14240534269101
$ws.Cells.Item(14, 18).Value = 667.281648084219
$ws.Cells.Item(14, 19).Value = 0.00961782610269485
$ws.Cells.Item(14, 20).Value = 0.009617826102694849

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 7.135311000000001
$ws.Cells.Item(15, 8).Value = 21.405933
$ws.Cells.Item(15, 9).Value = 0.115169939680053
$ws.Cells.Item(15, 10).Value = 0.115169939680053
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 35.71561933333334
$ws.Cells.Item(15, 14).Value = 107.146858
$ws.Cells.Item(15, 15).Value = 0.287039855156433
$ws.Cells.Item(15, 16).Value = 0.287039855156433
$ws.Cells.Item(15, 17).Value = 254.842051500946
$ws.Cells.Item(15, 18).Value = 2293.578463508514
$ws.Cells.Item(15, 19).Value = 0.03305836280413754
$ws.Cells.Item(15, 20).Value = 0.03305836280413753

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 7.135311000000001
$ws.Cells.Item(16, 8).Value = 21.405933
$ws.Cells.Item(16, 9).Value = 0.115169939680053
$ws.Cells.Item(16, 10).Value = 0.115169939680053
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 51.87044599999999
$ws.Cells.Item(16, 14).Value = 155.611338
$ws.Cells.Item(16, 15).Value = 0.4168732219867682
$ws.Cells.Item(16, 16).Value = 0.4168732219867682
$ws.Cells.Item(16, 17).Value = 370.111763918706
$ws.Cells.Item(16, 18).Value = 3331.005875268354
$ws.Cells.Item(16, 19).Value = 0.04801126383044543
$ws.Cells.Item(16, 20).Value = 0.04801126383044543

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 7.135311000000001
$ws.Cells.Item(17, 8).Value = 21.405933
$ws.Cells.Item(17, 9).Value = 0.115169939680053
$ws.Cells.Item(17, 10).Value = 0.115169939680053
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 26.450408
$ws.Cells.Item(17, 14).Value = 79.351224
$ws.Cells.Item(17, 15).Value = 0.2125770579613792
$ws.Cells.Item(17, 16).Value = 0.2125770579613792
$ws.Cells.Item(17, 17).Value = 188.731887156888
$ws.Cells.Item(17, 18).Value = 1698.586984411992
$ws.Cells.Item(17, 19).Value = 0.02448248694277517
$ws.Cells.Item(17, 20).Value = 0.02448248694277517
